$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.233.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.299.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.74%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.644"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.63"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +11.69%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  +13.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0986"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.41%  "

$ws.Range("E14").Value = "  +3.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.641.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.901"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.295.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.203.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.36%  "

$ws.Range("E20").Value = "  +7.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.46%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.58%  "

$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "

$ws.Range("E30").Value = "  +8.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.17%  "

$ws.Range("E34").Value = "  +8.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0815"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +27.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.126"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +20.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0314"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +23.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.38%  "

$ws.Range("E43").Value = "  +11.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.215"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +14.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "62.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.105"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.85%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.80%  "

$ws.Range("E50").Value = "  +0.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.14%  "
